# Update leve-profit computed columns (currentAveragePrice* / LevePrice* / LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed market data
# pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 500
$ws.Range("I48").Value = 500
$ws.Range("K48").Value = 1500
$ws.Range("M48").Value = -1208

$ws.Range("H56").Value = 500
$ws.Range("I56").Value = 500
$ws.Range("K56").Value = 1500
$ws.Range("M56").Value = -966

$ws.Range("H114").Value = 39719.668
$ws.Range("J114").Value = 39719.668
$ws.Range("L114").Value = 39719.668
$ws.Range("N114").Value = -48397.668

$ws.Range("H129").Value = 756.85
$ws.Range("I129").Value = 425.35
$ws.Range("K129").Value = 1276.05
$ws.Range("M129").Value = 3723.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20931.105
$ws.Range("I32").Value = 17544.025
$ws.Range("K32").Value = 17544.025
$ws.Range("M32").Value = -17257.025

$ws.Range("H44").Value = 20016.334
$ws.Range("J44").Value = 20016.334
$ws.Range("L44").Value = 20016.334
$ws.Range("N44").Value = -20992.334

$ws.Range("H45").Value = 1287.8462
$ws.Range("I45").Value = 1187.8572
$ws.Range("J45").Value = 1404.5
$ws.Range("K45").Value = 1187.8572
$ws.Range("L45").Value = 1404.5
$ws.Range("M45").Value = -810.8571999999999
$ws.Range("N45").Value = -2158.5

$ws.Range("H63").Value = 2622.25
$ws.Range("I63").Value = 2673.9285
$ws.Range("J63").Value = 2549.9
$ws.Range("K63").Value = 2673.9285
$ws.Range("L63").Value = 2549.9
$ws.Range("M63").Value = -1987.9285
$ws.Range("N63").Value = -3921.9

$ws.Range("H66").Value = 2622.25
$ws.Range("I66").Value = 2673.9285
$ws.Range("J66").Value = 2549.9
$ws.Range("K66").Value = 13369.6425
$ws.Range("L66").Value = 12749.5
$ws.Range("M66").Value = -9937.6425
$ws.Range("N66").Value = -19613.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13555.25
$ws.Range("I26").Value = 13555.25
$ws.Range("K26").Value = 13555.25
$ws.Range("M26").Value = -13263.25

$ws.Range("H107").Value = 2228.2222
$ws.Range("I107").Value = 1965.2632
$ws.Range("J107").Value = 2852.75
$ws.Range("K107").Value = 1965.2632
$ws.Range("L107").Value = 2852.75
$ws.Range("M107").Value = -45.2632000000001
$ws.Range("N107").Value = -6692.75

$ws.Range("H123").Value = 18569.215
$ws.Range("J123").Value = 18569.215
$ws.Range("L123").Value = 18569.215
$ws.Range("N123").Value = -28369.215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3978.673
$ws.Range("I31").Value = 1189.1212
$ws.Range("J31").Value = 8823.684999999999
$ws.Range("K31").Value = 1189.1212
$ws.Range("L31").Value = 8823.684999999999
$ws.Range("M31").Value = -894.1212
$ws.Range("N31").Value = -9413.684999999999

$ws.Range("H34").Value = 3978.673
$ws.Range("I34").Value = 1189.1212
$ws.Range("J34").Value = 8823.684999999999
$ws.Range("K34").Value = 1189.1212
$ws.Range("L34").Value = 8823.684999999999
$ws.Range("M34").Value = -987.1212
$ws.Range("N34").Value = -9227.684999999999

$ws.Range("H59").Value = 32323.143
$ws.Range("J59").Value = 32323.143
$ws.Range("L59").Value = 32323.143
$ws.Range("N59").Value = -34613.143

$ws.Range("H99").Value = 1131.3572
$ws.Range("I99").Value = 1011.5833
$ws.Range("J99").Value = 1850
$ws.Range("K99").Value = 1011.5833
$ws.Range("L99").Value = 1850
$ws.Range("M99").Value = 486.4167
$ws.Range("N99").Value = -4846

$ws.Range("H105").Value = 1955
$ws.Range("I105").Value = 1955
$ws.Range("K105").Value = 1955
$ws.Range("M105").Value = -208

$ws.Range("H126").Value = 1131.3572
$ws.Range("I126").Value = 1011.5833
$ws.Range("J126").Value = 1850
$ws.Range("K126").Value = 3034.7499
$ws.Range("L126").Value = 5550
$ws.Range("M126").Value = -564.7498999999998
$ws.Range("N126").Value = -10490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 920.1277
$ws.Range("J131").Value = 928.8
$ws.Range("L131").Value = 2786.4
$ws.Range("N131").Value = -12866.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2006
$ws.Range("I102").Value = 2012
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2012
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -390
$ws.Range("N102").Value = -5244

$ws.Range("H113").Value = 1661.1923
$ws.Range("I113").Value = 853.6
$ws.Range("J113").Value = 2762.4546
$ws.Range("K113").Value = 853.6
$ws.Range("L113").Value = 2762.4546
$ws.Range("M113").Value = 1316.4
$ws.Range("N113").Value = -7102.4546

$ws.Range("H126").Value = 1735.125
$ws.Range("I126").Value = 1666.75
$ws.Range("J126").Value = 1803.5
$ws.Range("K126").Value = 5000.25
$ws.Range("L126").Value = 5410.5
$ws.Range("M126").Value = -2530.25
$ws.Range("N126").Value = -10350.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8623169
$ws.Range("I7").Value = 11365966
$ws.Range("J7").Value = 2948.4285
$ws.Range("K7").Value = 11365966
$ws.Range("L7").Value = 2948.4285
$ws.Range("M7").Value = -11365854
$ws.Range("N7").Value = -3172.4285

$ws.Range("H40").Value = 3108.2778
$ws.Range("I40").Value = 3049.6155
$ws.Range("J40").Value = 3260.8
$ws.Range("K40").Value = 3049.6155
$ws.Range("L40").Value = 3260.8
$ws.Range("M40").Value = -2913.6155
$ws.Range("N40").Value = -3532.8

$ws.Range("H46").Value = 2525945.5
$ws.Range("I46").Value = 3367684
$ws.Range("J46").Value = 730
$ws.Range("K46").Value = 3367684
$ws.Range("L46").Value = 730
$ws.Range("M46").Value = -3367496
$ws.Range("N46").Value = -1106

$ws.Range("H55").Value = 311.1111
$ws.Range("I55").Value = 312.5
$ws.Range("J55").Value = 310
$ws.Range("K55").Value = 312.5
$ws.Range("L55").Value = 310
$ws.Range("M55").Value = -139.5
$ws.Range("N55").Value = -656

$ws.Range("H126").Value = 8623169
$ws.Range("I126").Value = 11365966
$ws.Range("J126").Value = 2948.4285
$ws.Range("K126").Value = 34097898
$ws.Range("L126").Value = 8845.2855
$ws.Range("M126").Value = -34095428
$ws.Range("N126").Value = -13785.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1896.7858
$ws.Range("I126").Value = 1888.8462
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5666.5386
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3196.5386
$ws.Range("N126").Value = -10940
